$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# FC Order ID (A2): purely numeric-looking text -> force text via TEXT()
# formula then flatten to a static value so it keeps its original cell
# style (and doesn't get reinterpreted as a Number).
$ws.Range("A2").Formula = "=TEXT(58572102,""0"")"
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)

# Tracking # (C2): not numeric-looking, plain text assignment is safe.
$ws.Range("C2").Value = "FCT943034220790415360"

# New Invoice Amount (F2): purely numeric-looking text -> same TEXT() trick.
$ws.Range("F2").Formula = "=TEXT(67.81,""0.00"")"
$ws.Range("F2").Copy()
$ws.Range("F2").PasteSpecial(-4163)

# SECONDARY INV # (I2): contains a "+", not parsed as a number, plain text is safe.
$ws.Range("I2").Value = "58572102+1"
